$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New measurement row of test data (row 16)
$ws.Range("A16").Value = 818
$ws.Range("B16").Value = 1016
$ws.Range("C16").Formula = "=B16-A16"
$ws.Range("D16").Value = 16
$ws.Range("E16").Formula = "=C16/D16"
$ws.Range("G16").Value = "good"

# The old placeholder row 17 (empty except for formatting) is removed
$ws.Range("G17").Clear()
